$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the "Case" rows in the order needed so that the shared-string table
# ends up in the same order as the target workbook.
$ws.Range("S19").Value = 'Code Owners,316049311212458244,  Case,9,  Total waiting time: 62249.34007431242,  average waiting time per call: 62.24934007431242,  unCompleted calls,6,  certificate, -316049312923083330'
$ws.Range("S18").Value = 'Code Owners,316049311212458244,  Case,8,  Total waiting time: 139783.88209694857,  average waiting time per call: 139.78388209694856,  unCompleted calls,23,  certificate, -316049315692285783'
$ws.Range("S17").Value = 'Code Owners,316049311212458244,  Case,7,  Total waiting time: 177690.12115705168,  average waiting time per call: 177.69012115705166,  unCompleted calls,37,  certificate, -316049316911369233'
$ws.Range("S16").Value = 'Code Owners,316049311212458244,  Case,6,  Total waiting time: 66957.88209694847,  average waiting time per call: 66.95788209694847,  unCompleted calls,15,  certificate, -316049308873405043'
$ws.Range("S15").Value = 'Code Owners,316049311212458244,  Case,5,  Total waiting time: 78504.1211570518,  average waiting time per call: 78.50412115705181,  unCompleted calls,22,  certificate, -316049309455165983'
$ws.Range("S14").Value = 'Code Owners,316049311212458244,  Case,4,  Total waiting time: 23197.455368642088,  average waiting time per call: 46.39491073728418,  unCompleted calls,2,  certificate, -316049312499393040'
$ws.Range("S13").Value = 'Code Owners,316049311212458244,  Case,3,  Total waiting time: 19723.538284333124,  average waiting time per call: 49.30884571083281,  unCompleted calls,2,  certificate, -316049312614093487'
$ws.Range("S12").Value = 'Code Owners,316049311212458244,  Case,2,  Total waiting time: 5051.792822120196,  average waiting time per call: 50.51792822120196,  unCompleted calls,6,  certificate, -316049312646921093'
$ws.Range("S10").Value = 'Code Owners,316049311212458244,  Case,0,  Total waiting time: 231.9897426188186,  average waiting time per call: 23.19897426188186,  unCompleted calls,0,  certificate, -316049311681347185'
$ws.Range("S11").Value = 'Code Owners,316049311212458244,  Case,1,  Total waiting time: 327.9897426188186,  average waiting time per call: 32.79897426188186,  unCompleted calls,4,  certificate, -316049312126102141'

# New header / submitter rows inserted above the case list
$ws.Range("Q8").Value = 'מגישים:'
$ws.Range("R8").Value = 316049311
$ws.Range("S8").Value = 'אמיר סבג'
$ws.Range("R9").Value = 212458244
$ws.Range("S9").Value = 'אורי דרשן'

# Update the active selection to match the target view
$null = $ws.Range("S9").Select()
